$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOT2045")

# ---------------------------------------------------------------------------
# 1) Update activation date 01/01/2018 -> 01/01/2022
#    NOTE: assigning a date-shaped literal straight to .Value auto-converts
#    it to a real date serial (and a new date number-format style), which is
#    not what the source workbook has (it stores it as a plain text string
#    sharing the original style). To keep it as literal text with the
#    original style untouched, build the text in a scratch cell via a
#    formula (so it is never auto-parsed as a date), paste only its
#    *value* into the target cell, then wipe the scratch cell.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="01/01/" & "2022"'
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# ---------------------------------------------------------------------------
# 2) Add the English translation of the "Objetivos:" answer into B11/C11
$objText = "Provide students with the knowledge of cell biology necessary to understand the other subjects of the course and the training of the Environmental Engineer."
$ws.Range("B11").Value = $objText
$ws.Range("C11").Value = $objText
# Match the formatting used by the analogous existing row (B10/C10)
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Update "Programa resumido:" text (drop origin/evolution clause)
$resumido = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# 4) Update "Short syllabus:" text (drop origin/evolution clause)
$shortSyllabus = "Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division."
$ws.Range("B15").Value = $shortSyllabus
$ws.Range("C15").Value = $shortSyllabus

# 5) Update "Programa:" text (replace origin/evolution bullet with new evolutionary-history bullet)
$programa = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# 6) Update "Syllabus:" text (replace origin/evolution clause with new evolutionary-history clause)
$syllabus = "Cell structure and evolutionary history: prokaryotic microorganisms andeukaryotic and their evolutionary relationships between the Bacteria, Archaea andEukarya.Microscope analysis of cells structure: optical and electron microscope.Structure and function of major organic molecules: carbohydrates, lipids, nucleic acids and proteins. Internal organization of the cell: membrane structure and function; intracelular compartments and protein sorting; vesicular traffic (endocytosis and exocytosis).Nucleus and genetic material organization: structure and functionCell cycle and cell division: mitosis and meiosisCell energy conversion: mitochondria and chloroplast."
$ws.Range("B17").Value = $syllabus
$ws.Range("C17").Value = $syllabus

$wb.Save()
